$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column C ("max") entirely, shifting D ("prediction") and E ("rejection-f") left
$ws.Range("C:C").Delete()

# Update the B column values (was all 1, now actual max-score-like numbers)
$ws.Range("B2").Value = 1643.945200406717
$ws.Range("B3").Value = 1535.092286422305
$ws.Range("B4").Value = 1723.472007440319
